# Updates cryptos price/volume data (Price column D, Volume(1h) column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value would otherwise be auto-coerced to a number
# by Excel's input parsing (single-dot decimal-looking strings); format as Text first
# so the literal string is preserved, matching the source data's inline-string type.
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D17", "D18", "D19", "D21", "D22", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Column D (Price) updates
$ws.Range("D2").Value = "28.281.70"
$ws.Range("D3").Value = "1.816.81"
$ws.Range("D5").Value = "326.10"
$ws.Range("D6").Value = "1.0000"
$ws.Range("D7").Value = "0.4367"
$ws.Range("D8").Value = "0.3672"
$ws.Range("D9").Value = "44.88"
$ws.Range("D10").Value = "0.07681"
$ws.Range("D11").Value = "1.142"
$ws.Range("D12").Value = "0.9999"
$ws.Range("D13").Value = "22.06"
$ws.Range("D14").Value = "6.319"
$ws.Range("D15").Value = "7.501"
$ws.Range("D16").Value = "1.819.56"
$ws.Range("D17").Value = "95.12"
$ws.Range("D18").Value = "0.00001081"
$ws.Range("D19").Value = "0.06497"
$ws.Range("D21").Value = "17.40"
$ws.Range("D22").Value = "6.242"
$ws.Range("D23").Value = "28.290.02"
$ws.Range("D24").Value = "11.58"
$ws.Range("D25").Value = "2.130"
$ws.Range("D26").Value = "161.12"
$ws.Range("D27").Value = "20.73"
$ws.Range("D28").Value = "2.025.90"
$ws.Range("D29").Value = "2.285"
$ws.Range("D30").Value = "130.21"
$ws.Range("D31").Value = "1.215"
$ws.Range("D32").Value = "6.012"
$ws.Range("D33").Value = "0.09147"
$ws.Range("D34").Value = "3.550"
$ws.Range("D35").Value = "13.07"
$ws.Range("D36").Value = "0.02374"
$ws.Range("D37").Value = "5.253"
$ws.Range("D38").Value = "0.2178"
$ws.Range("D39").Value = "0.6608"
$ws.Range("D40").Value = "0.06213"
$ws.Range("D41").Value = "1.200"
$ws.Range("D42").Value = "8.080"
$ws.Range("D43").Value = "1.433"
$ws.Range("D44").Value = "0.9988"
$ws.Range("D45").Value = "13.87"
$ws.Range("D46").Value = "0.6124"
$ws.Range("D47").Value = "3.738"
$ws.Range("D48").Value = "2.023"
$ws.Range("D49").Value = "125.74"
$ws.Range("D50").Value = "1.165"
$ws.Range("D51").Value = "0.06997"

# Column E (Volume(1h)) updates
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("E3").Value = "  +3.67%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  +1.22%  "
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("E10").Value = "  +2.97%  "
$ws.Range("E11").Value = "  +1.98%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("E13").Value = "  +2.23%  "
$ws.Range("E14").Value = "  +2.66%  "
$ws.Range("E15").Value = "  +3.59%  "
$ws.Range("E16").Value = "  +3.84%  "
$ws.Range("E17").Value = "  +8.18%  "
$ws.Range("E18").Value = "  +1.42%  "
$ws.Range("E19").Value = "  +4.84%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("E23").Value = "  +2.14%  "
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("E25").Value = "  -8.54%  "
$ws.Range("E26").Value = "  +4.88%  "
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("E28").Value = "  +3.84%  "
$ws.Range("E29").Value = "  -3.25%  "
$ws.Range("E30").Value = "  +2.31%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("E32").Value = "  +5.16%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("E34").Value = "  -1.91%  "
$ws.Range("E35").Value = "  +3.55%  "
$ws.Range("E36").Value = "  +2.84%  "
$ws.Range("E37").Value = "  +2.95%  "
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("E39").Value = "  +2.12%  "
$ws.Range("E40").Value = "  +1.97%  "
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("E42").Value = "  +1.63%  "
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("E46").Value = "  +3.21%  "
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("E48").Value = "  +2.48%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  +4.06%  "
$ws.Range("E51").Value = "  +1.44%  "
